$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 218-222 with revised figures (columns B..Q); R is unchanged.
$updates = @{
    218 = @(3782, 1205, 895, 316, -6, 6372, 2991, 3381, 2251, 1130, -1395, -2289, 702, -141, -2841, -9)
    219 = @(5414, 6189, 5790, 469, -71, -402, 436, -838, -934, 96, -1325, 266, 303, 17, -45, -9)
    220 = @(1851, 930, 563, 196, 171, -738, -140, -598, -115, -483, -970, 1749, -35, 86, 1706, -9)
    221 = @(4040, 537, 122, 418, -3, -625, 1019, -1644, 513, -2156, -1086, 3105, 1212, 6, 1895, -9)
    222 = @(247, 368, 111, 77, 180, -5217, -6367, 1151, 957, 194, -789, 953, 531, 101, 228, 94)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($c = 2; $c -le 17; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 2]
    }
}

# Append new row 223 with the 01-06-2021 series.
$newRow = 223
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "01-06-2021"
$ws.Cells.Item($newRow, 1).Style = "Normal"
$rowValues = @(-1796, 544, 234, 290, 20, 486, 2261, -1776, -502, -1273, -856, 626, -445, -29, 1109, -9, -2596)
for ($c = 2; $c -le 18; $c++) {
    $ws.Cells.Item($newRow, $c).Value = $rowValues[$c - 2]
}
